# Auto-generated data refresh for Golem_Profits workbook
# Applies updated market-price derived values (columns H-N) per the scheduled runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 126.09091  # ALC!H6  was 115.583336
$ws.Cells.Item(6, 9).Value = 46.375  # ALC!I6  was 38
$ws.Cells.Item(6, 10).Value = 338.66666  # ALC!J6  was 503.5
$ws.Cells.Item(6, 11).Value = 139.125  # ALC!K6  was 114
$ws.Cells.Item(6, 12).Value = 1015.99998  # ALC!L6  was 1510.5
$ws.Cells.Item(6, 13).Value = -27.125  # ALC!M6  was -2
$ws.Cells.Item(6, 14).Value = -1239.99998  # ALC!N6  was -1734.5
$ws.Cells.Item(9, 8).Value = 1089.1666  # ALC!H9  was 843
$ws.Cells.Item(9, 9).Value = 299  # ALC!I9  was 309.8
$ws.Cells.Item(9, 10).Value = 1247.2  # ALC!J9  was 1509.5
$ws.Cells.Item(9, 11).Value = 299  # ALC!K9  was 309.8
$ws.Cells.Item(9, 12).Value = 1247.2  # ALC!L9  was 1509.5
$ws.Cells.Item(9, 13).Value = -130  # ALC!M9  was -140.8
$ws.Cells.Item(9, 14).Value = -1585.2  # ALC!N9  was -1847.5
$ws.Cells.Item(12, 8).Value = 969.7  # ALC!H12  was 932.9167
$ws.Cells.Item(12, 9).Value = 956.7143  # ALC!I12  was 910.55554
$ws.Cells.Item(12, 11).Value = 956.7143  # ALC!K12  was 910.55554
$ws.Cells.Item(12, 13).Value = -786.7143  # ALC!M12  was -740.55554
$ws.Cells.Item(17, 8).Value = 3066.3333  # ALC!H17  was 3374.25
$ws.Cells.Item(17, 10).Value = 3074.625  # ALC!J17  was 3499
$ws.Cells.Item(17, 12).Value = 9223.875  # ALC!L17  was 10497
$ws.Cells.Item(17, 14).Value = -9559.875  # ALC!N17  was -10833
$ws.Cells.Item(21, 8).Value = 2139  # ALC!H21  was 417
$ws.Cells.Item(21, 10).Value = 3000  # ALC!J21  was 0
$ws.Cells.Item(21, 12).Value = 3000  # ALC!L21  was 0
$ws.Cells.Item(21, 14).Value = -3936  # ALC!N21  was None
$ws.Cells.Item(23, 8).Value = 2139  # ALC!H23  was 417
$ws.Cells.Item(23, 10).Value = 3000  # ALC!J23  was 0
$ws.Cells.Item(23, 12).Value = 3000  # ALC!L23  was 0
$ws.Cells.Item(23, 14).Value = -3468  # ALC!N23  was None
$ws.Cells.Item(33, 8).Value = 494.07144  # ALC!H33  was 413.17648
$ws.Cells.Item(33, 9).Value = 455.30768  # ALC!I33  was 376.625
$ws.Cells.Item(33, 11).Value = 455.30768  # ALC!K33  was 376.625
$ws.Cells.Item(33, 13).Value = -226.30768  # ALC!M33  was -147.625
$ws.Cells.Item(40, 8).Value = 2900  # ALC!H40  was 9450.5
$ws.Cells.Item(40, 9).Value = 2900  # ALC!I40  was 9450.5
$ws.Cells.Item(40, 11).Value = 2900  # ALC!K40  was 9450.5
$ws.Cells.Item(40, 13).Value = -2725  # ALC!M40  was -9275.5
$ws.Cells.Item(70, 8).Value = 3500  # ALC!H70  was 3333.3333
$ws.Cells.Item(70, 10).Value = 0  # ALC!J70  was 2500
$ws.Cells.Item(70, 12).Value = 0  # ALC!L70  was 7500
$ws.Cells.Item(70, 14).Value = $null  # ALC!N70  was -8040
$ws.Cells.Item(73, 8).Value = 3500  # ALC!H73  was 3333.3333
$ws.Cells.Item(73, 10).Value = 0  # ALC!J73  was 2500
$ws.Cells.Item(73, 12).Value = 0  # ALC!L73  was 7500
$ws.Cells.Item(73, 14).Value = $null  # ALC!N73  was -9372
$ws.Cells.Item(80, 8).Value = 498.33334  # ALC!H80  was 647.5
$ws.Cells.Item(80, 10).Value = 498.33334  # ALC!J80  was 647.5
$ws.Cells.Item(80, 12).Value = 1495.00002  # ALC!L80  was 1942.5
$ws.Cells.Item(80, 14).Value = -3491.00002  # ALC!N80  was -3938.5
$ws.Cells.Item(83, 8).Value = 498.33334  # ALC!H83  was 647.5
$ws.Cells.Item(83, 10).Value = 498.33334  # ALC!J83  was 647.5
$ws.Cells.Item(83, 12).Value = 4485.00006  # ALC!L83  was 5827.5
$ws.Cells.Item(83, 14).Value = -14469.00006  # ALC!N83  was -15811.5
$ws.Cells.Item(113, 8).Value = 3938.8  # ALC!H113  was 3948.75
$ws.Cells.Item(113, 10).Value = 3899.5  # ALC!J113  was 3900
$ws.Cells.Item(113, 12).Value = 3899.5  # ALC!L113  was 3900
$ws.Cells.Item(113, 14).Value = -10407.5  # ALC!N113  was -10408
$ws.Cells.Item(116, 8).Value = 4200  # ALC!H116  was 4376.25
$ws.Cells.Item(116, 9).Value = 4200  # ALC!I116  was 4376.25
$ws.Cells.Item(116, 11).Value = 4200  # ALC!K116  was 4376.25
$ws.Cells.Item(116, 13).Value = -758  # ALC!M116  was -934.25
$ws.Cells.Item(132, 8).Value = 1942.4286  # ALC!H132  was 1669.8
$ws.Cells.Item(132, 9).Value = 1516.1666  # ALC!I132  was 1669.8
$ws.Cells.Item(132, 10).Value = 4500  # ALC!J132  was 0
$ws.Cells.Item(132, 11).Value = 4548.4998  # ALC!K132  was 5009.4
$ws.Cells.Item(132, 12).Value = 13500  # ALC!L132  was 0
$ws.Cells.Item(132, 13).Value = -2018.4998  # ALC!M132  was -2479.4
$ws.Cells.Item(132, 14).Value = -18560  # ALC!N132  was None

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1995.6  # ARM!H2  was 1999.5
$ws.Cells.Item(2, 9).Value = 1995.6  # ARM!I2  was 1999.5
$ws.Cells.Item(2, 11).Value = 1995.6  # ARM!K2  was 1999.5
$ws.Cells.Item(2, 13).Value = -1882.6  # ARM!M2  was -1886.5
$ws.Cells.Item(63, 8).Value = 5343.3335  # ARM!H63  was 1750
$ws.Cells.Item(63, 9).Value = 2410.8  # ARM!I63  was 1750
$ws.Cells.Item(63, 10).Value = 20006  # ARM!J63  was 0
$ws.Cells.Item(63, 11).Value = 2410.8  # ARM!K63  was 1750
$ws.Cells.Item(63, 12).Value = 20006  # ARM!L63  was 0
$ws.Cells.Item(63, 13).Value = -1724.8  # ARM!M63  was -1064
$ws.Cells.Item(63, 14).Value = -21378  # ARM!N63  was None
$ws.Cells.Item(66, 8).Value = 5343.3335  # ARM!H66  was 1750
$ws.Cells.Item(66, 9).Value = 2410.8  # ARM!I66  was 1750
$ws.Cells.Item(66, 10).Value = 20006  # ARM!J66  was 0
$ws.Cells.Item(66, 11).Value = 12054  # ARM!K66  was 8750
$ws.Cells.Item(66, 12).Value = 100030  # ARM!L66  was 0
$ws.Cells.Item(66, 13).Value = -8622  # ARM!M66  was -5318
$ws.Cells.Item(66, 14).Value = -106894  # ARM!N66  was None
$ws.Cells.Item(116, 8).Value = 1995.6  # ARM!H116  was 1999.5
$ws.Cells.Item(116, 9).Value = 1995.6  # ARM!I116  was 1999.5
$ws.Cells.Item(116, 11).Value = 1995.6  # ARM!K116  was 1999.5
$ws.Cells.Item(116, 13).Value = 298.4000000000001  # ARM!M116  was 294.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1995.6  # BSM!H3  was 1999.5
$ws.Cells.Item(3, 9).Value = 1995.6  # BSM!I3  was 1999.5
$ws.Cells.Item(3, 11).Value = 1995.6  # BSM!K3  was 1999.5
$ws.Cells.Item(3, 13).Value = -1881.6  # BSM!M3  was -1885.5
$ws.Cells.Item(29, 8).Value = 984.9  # BSM!H29  was 987.4167

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 1599.6666  # CRP!H58  was 1666.1666
$ws.Cells.Item(58, 9).Value = 900  # CRP!I58  was 1000
$ws.Cells.Item(58, 10).Value = 2999  # CRP!J58  was 1999.25
$ws.Cells.Item(58, 11).Value = 900  # CRP!K58  was 1000
$ws.Cells.Item(58, 12).Value = 2999  # CRP!L58  was 1999.25
$ws.Cells.Item(58, 13).Value = -697  # CRP!M58  was -797
$ws.Cells.Item(58, 14).Value = -3405  # CRP!N58  was -2405.25
$ws.Cells.Item(62, 8).Value = 3583.5  # CRP!H62  was 3845.875
$ws.Cells.Item(62, 9).Value = 3378  # CRP!I62  was 3693.6
$ws.Cells.Item(62, 10).Value = 4200  # CRP!J62  was 4099.6665
$ws.Cells.Item(62, 11).Value = 3378  # CRP!K62  was 3693.6
$ws.Cells.Item(62, 12).Value = 4200  # CRP!L62  was 4099.6665
$ws.Cells.Item(62, 13).Value = -2754  # CRP!M62  was -3069.6
$ws.Cells.Item(62, 14).Value = -5448  # CRP!N62  was -5347.6665
$ws.Cells.Item(65, 8).Value = 3583.5  # CRP!H65  was 3845.875
$ws.Cells.Item(65, 9).Value = 3378  # CRP!I65  was 3693.6
$ws.Cells.Item(65, 10).Value = 4200  # CRP!J65  was 4099.6665
$ws.Cells.Item(65, 11).Value = 16890  # CRP!K65  was 18468
$ws.Cells.Item(65, 12).Value = 21000  # CRP!L65  was 20498.3325
$ws.Cells.Item(65, 13).Value = -13770  # CRP!M65  was -15348
$ws.Cells.Item(65, 14).Value = -27240  # CRP!N65  was -26738.3325
$ws.Cells.Item(75, 8).Value = 65500  # CRP!H75  was 61500
$ws.Cells.Item(75, 10).Value = 65500  # CRP!J75  was 61500
$ws.Cells.Item(75, 12).Value = 65500  # CRP!L75  was 61500
$ws.Cells.Item(75, 14).Value = -67496  # CRP!N75  was -63496
$ws.Cells.Item(78, 8).Value = 65500  # CRP!H78  was 61500
$ws.Cells.Item(78, 10).Value = 65500  # CRP!J78  was 61500
$ws.Cells.Item(78, 12).Value = 196500  # CRP!L78  was 184500
$ws.Cells.Item(78, 14).Value = -206484  # CRP!N78  was -194484
$ws.Cells.Item(107, 8).Value = 1042.1428  # CRP!H107  was 909.4
$ws.Cells.Item(107, 9).Value = 850  # CRP!I107  was 742.7143
$ws.Cells.Item(107, 11).Value = 850  # CRP!K107  was 742.7143
$ws.Cells.Item(107, 13).Value = 1070  # CRP!M107  was 1177.2857
$ws.Cells.Item(136, 8).Value = 1599.6666  # CRP!H136  was 1666.1666
$ws.Cells.Item(136, 9).Value = 900  # CRP!I136  was 1000
$ws.Cells.Item(136, 10).Value = 2999  # CRP!J136  was 1999.25
$ws.Cells.Item(136, 11).Value = 2700  # CRP!K136  was 3000
$ws.Cells.Item(136, 12).Value = 8997  # CRP!L136  was 5997.75
$ws.Cells.Item(136, 13).Value = -150  # CRP!M136  was -450
$ws.Cells.Item(136, 14).Value = -14097  # CRP!N136  was -11097.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 154.7037  # CUL!H2  was 160.2963
$ws.Cells.Item(2, 10).Value = 234.27272  # CUL!J2  was 248
$ws.Cells.Item(2, 12).Value = 1405.63632  # CUL!L2  was 1488
$ws.Cells.Item(2, 14).Value = -1631.63632  # CUL!N2  was -1714
$ws.Cells.Item(23, 8).Value = 114  # CUL!H23  was 97.9375
$ws.Cells.Item(23, 9).Value = 78.5  # CUL!I23  was 59
$ws.Cells.Item(23, 10).Value = 129.77777  # CUL!J23  was 121.3
$ws.Cells.Item(23, 11).Value = 235.5  # CUL!K23  was 177
$ws.Cells.Item(23, 12).Value = 389.33331  # CUL!L23  was 363.9
$ws.Cells.Item(23, 13).Value = -0.5  # CUL!M23  was 58
$ws.Cells.Item(23, 14).Value = -859.33331  # CUL!N23  was -833.9
$ws.Cells.Item(25, 8).Value = 50  # CUL!H25  was 290
$ws.Cells.Item(25, 9).Value = 0  # CUL!I25  was 80
$ws.Cells.Item(25, 10).Value = 50  # CUL!J25  was 500
$ws.Cells.Item(25, 11).Value = 0  # CUL!K25  was 240
$ws.Cells.Item(25, 12).Value = 150  # CUL!L25  was 1500
$ws.Cells.Item(25, 13).Value = $null  # CUL!M25  was -71
$ws.Cells.Item(25, 14).Value = -488  # CUL!N25  was -1838
$ws.Cells.Item(26, 8).Value = 375  # CUL!H26  was 1043.75
$ws.Cells.Item(26, 9).Value = 375  # CUL!I26  was 391.66666
$ws.Cells.Item(26, 10).Value = 0  # CUL!J26  was 3000
$ws.Cells.Item(26, 11).Value = 1125  # CUL!K26  was 1174.99998
$ws.Cells.Item(26, 12).Value = 0  # CUL!L26  was 9000
$ws.Cells.Item(26, 13).Value = -837  # CUL!M26  was -886.9999800000001
$ws.Cells.Item(26, 14).Value = $null  # CUL!N26  was -9576
$ws.Cells.Item(30, 8).Value = 50  # CUL!H30  was 290
$ws.Cells.Item(30, 9).Value = 0  # CUL!I30  was 80
$ws.Cells.Item(30, 10).Value = 50  # CUL!J30  was 500
$ws.Cells.Item(30, 11).Value = 0  # CUL!K30  was 240
$ws.Cells.Item(30, 12).Value = 150  # CUL!L30  was 1500
$ws.Cells.Item(30, 13).Value = $null  # CUL!M30  was -138
$ws.Cells.Item(30, 14).Value = -354  # CUL!N30  was -1704
$ws.Cells.Item(38, 8).Value = 133.5  # CUL!H38  was 50.166668
$ws.Cells.Item(38, 9).Value = 45  # CUL!I38  was 62.4
$ws.Cells.Item(38, 10).Value = 222  # CUL!J38  was 41.42857
$ws.Cells.Item(38, 11).Value = 135  # CUL!K38  was 187.2
$ws.Cells.Item(38, 12).Value = 666  # CUL!L38  was 124.28571
$ws.Cells.Item(38, 13).Value = 212  # CUL!M38  was 159.8
$ws.Cells.Item(38, 14).Value = -1360  # CUL!N38  was -818.28571
$ws.Cells.Item(131, 8).Value = 1900  # CUL!H131  was 1840.2
$ws.Cells.Item(131, 9).Value = 0  # CUL!I131  was 1501
$ws.Cells.Item(131, 10).Value = 1900  # CUL!J131  was 1925
$ws.Cells.Item(131, 11).Value = 0  # CUL!K131  was 4503
$ws.Cells.Item(131, 12).Value = 5700  # CUL!L131  was 5775
$ws.Cells.Item(131, 13).Value = $null  # CUL!M131  was 537
$ws.Cells.Item(131, 14).Value = -15780  # CUL!N131  was -15855
$ws.Cells.Item(140, 8).Value = 2144.2  # CUL!H140  was 2430.25
$ws.Cells.Item(140, 9).Value = 2144.2  # CUL!I140  was 2430.25
$ws.Cells.Item(140, 11).Value = 6432.599999999999  # CUL!K140  was 7290.75
$ws.Cells.Item(140, 13).Value = -1252.599999999999  # CUL!M140  was -2110.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(29, 8).Value = 0  # GSM!H29  was 500
$ws.Cells.Item(29, 9).Value = 0  # GSM!I29  was 500
$ws.Cells.Item(29, 11).Value = 0  # GSM!K29  was 500
$ws.Cells.Item(29, 13).Value = $null  # GSM!M29  was -210
$ws.Cells.Item(107, 8).Value = 870.4545000000001  # GSM!H107  was 869.2727
$ws.Cells.Item(107, 9).Value = 162.83333  # GSM!I107  was 152
$ws.Cells.Item(107, 10).Value = 1719.6  # GSM!J107  was 2124.5
$ws.Cells.Item(107, 11).Value = 162.83333  # GSM!K107  was 152
$ws.Cells.Item(107, 12).Value = 1719.6  # GSM!L107  was 2124.5
$ws.Cells.Item(107, 13).Value = 1757.16667  # GSM!M107  was 1768
$ws.Cells.Item(107, 14).Value = -5559.6  # GSM!N107  was -5964.5
$ws.Cells.Item(132, 8).Value = 4014  # GSM!H132  was 1777.4
$ws.Cells.Item(132, 9).Value = 0  # GSM!I132  was 1218.25
$ws.Cells.Item(132, 11).Value = 0  # GSM!K132  was 3654.75
$ws.Cells.Item(132, 13).Value = $null  # GSM!M132  was -1124.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 920.44446  # LTW!H46  was 1165.1666
$ws.Cells.Item(46, 9).Value = 861.8182  # LTW!I46  was 998.2
$ws.Cells.Item(46, 10).Value = 1012.5714  # LTW!J46  was 2000
$ws.Cells.Item(46, 11).Value = 861.8182  # LTW!K46  was 998.2
$ws.Cells.Item(46, 12).Value = 1012.5714  # LTW!L46  was 2000
$ws.Cells.Item(46, 13).Value = -673.8182  # LTW!M46  was -810.2
$ws.Cells.Item(46, 14).Value = -1388.5714  # LTW!N46  was -2376
$ws.Cells.Item(55, 8).Value = 533.3333  # LTW!H55  was 503.72726
$ws.Cells.Item(55, 9).Value = 335.1111  # LTW!I55  was 350.66666
$ws.Cells.Item(55, 10).Value = 1128  # LTW!J55  was 1192.5
$ws.Cells.Item(55, 11).Value = 335.1111  # LTW!K55  was 350.66666
$ws.Cells.Item(55, 12).Value = 1128  # LTW!L55  was 1192.5
$ws.Cells.Item(55, 13).Value = -162.1111  # LTW!M55  was -177.66666
$ws.Cells.Item(55, 14).Value = -1474  # LTW!N55  was -1538.5
$ws.Cells.Item(68, 8).Value = 6000  # LTW!H68  was 5999.6665
$ws.Cells.Item(68, 9).Value = 6000  # LTW!I68  was 5999.6665
$ws.Cells.Item(68, 11).Value = 6000  # LTW!K68  was 5999.6665
$ws.Cells.Item(68, 13).Value = -5251  # LTW!M68  was -5250.6665
$ws.Cells.Item(71, 8).Value = 6000  # LTW!H71  was 5999.6665
$ws.Cells.Item(71, 9).Value = 6000  # LTW!I71  was 5999.6665
$ws.Cells.Item(71, 11).Value = 30000  # LTW!K71  was 29998.3325
$ws.Cells.Item(71, 13).Value = -26256  # LTW!M71  was -26254.3325
$ws.Cells.Item(109, 8).Value = 0  # LTW!H109  was 1
$ws.Cells.Item(109, 10).Value = 0  # LTW!J109  was 1
$ws.Cells.Item(109, 12).Value = 0  # LTW!L109  was 1
$ws.Cells.Item(109, 14).Value = $null  # LTW!N109  was -2775
$ws.Cells.Item(122, 8).Value = 3906  # LTW!H122  was 3921.923
$ws.Cells.Item(122, 9).Value = 3748.5  # LTW!I122  was 3765
$ws.Cells.Item(122, 11).Value = 11245.5  # LTW!K122  was 11295
$ws.Cells.Item(122, 13).Value = -8795.5  # LTW!M122  was -8845

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 0  # WVR!H62  was 5000
$ws.Cells.Item(62, 10).Value = 0  # WVR!J62  was 5000
$ws.Cells.Item(62, 12).Value = 0  # WVR!L62  was 5000
$ws.Cells.Item(62, 14).Value = $null  # WVR!N62  was -6248
$ws.Cells.Item(65, 8).Value = 0  # WVR!H65  was 5000
$ws.Cells.Item(65, 10).Value = 0  # WVR!J65  was 5000
$ws.Cells.Item(65, 12).Value = 0  # WVR!L65  was 25000
$ws.Cells.Item(65, 14).Value = $null  # WVR!N65  was -31240
$ws.Cells.Item(126, 8).Value = 4144.3  # WVR!H126  was 3663.7144
$ws.Cells.Item(126, 9).Value = 2974.3333  # WVR!I126  was 2711.111
$ws.Cells.Item(126, 10).Value = 5899.25  # WVR!J126  was 5378.4
$ws.Cells.Item(126, 11).Value = 8922.999899999999  # WVR!K126  was 8133.333
$ws.Cells.Item(126, 12).Value = 17697.75  # WVR!L126  was 16135.2
$ws.Cells.Item(126, 13).Value = -6452.999899999999  # WVR!M126  was -5663.333
$ws.Cells.Item(126, 14).Value = -22637.75  # WVR!N126  was -21075.2
